# chore: update Sheets via scheduled runner
# Refresh currentAveragePrice / LevePrice / LeveProfit market data columns
# (H, I, J, K, L, M, N) across the per-job profit sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 494.35715
$ws.Range("I33").Value = 494.35715
$ws.Range("K33").Value = 494.35715
$ws.Range("M33").Value = -265.35715

$ws.Range("H76").Value = 3070.9678
$ws.Range("I76").Value = 3084.6155
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 3084.6155
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -2769.6155
$ws.Range("N76").Value = -3630

$ws.Range("H79").Value = 3070.9678
$ws.Range("I79").Value = 3084.6155
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 3084.6155
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -1992.6155
$ws.Range("N79").Value = -5184

$ws.Range("H88").Value = 76125
$ws.Range("J88").Value = 101000
$ws.Range("L88").Value = 101000
$ws.Range("N88").Value = -101812

$ws.Range("H91").Value = 76125
$ws.Range("J91").Value = 101000
$ws.Range("L91").Value = 101000
$ws.Range("N91").Value = -103808

$ws.Range("H135").Value = 1300.8948
$ws.Range("I135").Value = 517.5333000000001
$ws.Range("J135").Value = 4238.5
$ws.Range("K135").Value = 4657.7997
$ws.Range("L135").Value = 38146.5
$ws.Range("M135").Value = -2122.7997
$ws.Range("N135").Value = -43216.5

$ws.Range("H137").Value = 2816.6182
$ws.Range("I137").Value = 2637.5122
$ws.Range("K137").Value = 7912.5366
$ws.Range("M137").Value = -5362.5366

$ws.Range("H138").Value = 2414.4722
$ws.Range("J138").Value = 3864.7097
$ws.Range("L138").Value = 11594.1291
$ws.Range("N138").Value = -21874.1291

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 440106.06
$ws.Range("I32").Value = 468853.84
$ws.Range("J32").Value = 28054.666
$ws.Range("K32").Value = 468853.84
$ws.Range("L32").Value = 28054.666
$ws.Range("M32").Value = -468566.84
$ws.Range("N32").Value = -28628.666

$ws.Range("H61").Value = 23812512
$ws.Range("I61").Value = 47621516
$ws.Range("J61").Value = 3507.1904
$ws.Range("K61").Value = 47621516
$ws.Range("L61").Value = 3507.1904
$ws.Range("M61").Value = -47621304
$ws.Range("N61").Value = -3931.1904

$ws.Range("H88").Value = 3555
$ws.Range("J88").Value = 3555
$ws.Range("L88").Value = 3555
$ws.Range("N88").Value = -4367

$ws.Range("H91").Value = 3555
$ws.Range("J91").Value = 3555
$ws.Range("L91").Value = 3555
$ws.Range("N91").Value = -6363

$ws.Range("H97").Value = 1276.8182
$ws.Range("I97").Value = 1140.2941
$ws.Range("J97").Value = 1741
$ws.Range("K97").Value = 1140.2941
$ws.Range("L97").Value = 1741
$ws.Range("M97").Value = -644.2941000000001
$ws.Range("N97").Value = -2733

$ws.Range("H113").Value = 44296
$ws.Range("J113").Value = 44296
$ws.Range("L113").Value = 44296
$ws.Range("N113").Value = -52974

$ws.Range("H132").Value = 3052.3962
$ws.Range("I132").Value = 2145.3142
$ws.Range("J132").Value = 4816.1665
$ws.Range("K132").Value = 6435.942599999999
$ws.Range("L132").Value = 14448.4995
$ws.Range("M132").Value = -3905.942599999999
$ws.Range("N132").Value = -19508.4995

$ws.Range("H136").Value = 23812512
$ws.Range("I136").Value = 47621516
$ws.Range("J136").Value = 3507.1904
$ws.Range("K136").Value = 142864548
$ws.Range("L136").Value = 10521.5712
$ws.Range("M136").Value = -142861998
$ws.Range("N136").Value = -15621.5712

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2092.691
$ws.Range("I134").Value = 1494.1282
$ws.Range("J134").Value = 3551.6875
$ws.Range("K134").Value = 4482.3846
$ws.Range("L134").Value = 10655.0625
$ws.Range("M134").Value = -1947.3846
$ws.Range("N134").Value = -15725.0625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14692.315
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 14692.315
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 14692.315
$ws.Range("N31").Value = -15282.315
$ws.Range("M31").ClearContents()

$ws.Range("H34").Value = 14692.315
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 14692.315
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 14692.315
$ws.Range("N34").Value = -15096.315
$ws.Range("M34").ClearContents()

$ws.Range("H58").Value = 1156.1951
$ws.Range("I58").Value = 816.4828
$ws.Range("J58").Value = 1977.1666
$ws.Range("K58").Value = 816.4828
$ws.Range("L58").Value = 1977.1666
$ws.Range("M58").Value = -613.4828
$ws.Range("N58").Value = -2383.1666

$ws.Range("H132").Value = 18230268
$ws.Range("I132").Value = 25000998
$ws.Range("K132").Value = 75002994
$ws.Range("M132").Value = -75000464

$ws.Range("H136").Value = 1156.1951
$ws.Range("I136").Value = 816.4828
$ws.Range("J136").Value = 1977.1666
$ws.Range("K136").Value = 2449.4484
$ws.Range("L136").Value = 5931.4998
$ws.Range("M136").Value = 100.5515999999998
$ws.Range("N136").Value = -11031.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1885.2
$ws.Range("J5").Value = 2750
$ws.Range("L5").Value = 8250
$ws.Range("N5").Value = -8474

$ws.Range("H68").Value = 1022.55554
$ws.Range("I68").Value = 1200
$ws.Range("J68").Value = 971.8570999999999
$ws.Range("K68").Value = 3600
$ws.Range("L68").Value = 2915.5713
$ws.Range("M68").Value = -2789
$ws.Range("N68").Value = -4537.5713

$ws.Range("H71").Value = 1022.55554
$ws.Range("I71").Value = 1200
$ws.Range("J71").Value = 971.8570999999999
$ws.Range("K71").Value = 10800
$ws.Range("L71").Value = 8746.713899999999
$ws.Range("M71").Value = -6744
$ws.Range("N71").Value = -16858.7139

$ws.Range("H76").Value = 3999.7104
$ws.Range("I76").Value = 3989
$ws.Range("K76").Value = 11967
$ws.Range("M76").Value = -11584

$ws.Range("H79").Value = 3999.7104
$ws.Range("I79").Value = 3989
$ws.Range("K79").Value = 11967
$ws.Range("M79").Value = -10641

$ws.Range("H122").Value = 2816.111
$ws.Range("I122").Value = 348.52
$ws.Range("J122").Value = 5900.6
$ws.Range("K122").Value = 3136.68
$ws.Range("L122").Value = 53105.4
$ws.Range("M122").Value = -686.6799999999998
$ws.Range("N122").Value = -58005.4

$ws.Range("H129").Value = 1220.5652
$ws.Range("I129").Value = 476.66666
$ws.Range("J129").Value = 1483.1177
$ws.Range("K129").Value = 1429.99998
$ws.Range("L129").Value = 4449.3531
$ws.Range("M129").Value = 3570.00002
$ws.Range("N129").Value = -14449.3531

$ws.Range("H135").Value = 1885.2
$ws.Range("J135").Value = 2750
$ws.Range("L135").Value = 24750
$ws.Range("N135").Value = -29820

$ws.Range("H137").Value = 11915236
$ws.Range("I137").Value = 55589224
$ws.Range("J137").Value = 4148.4546
$ws.Range("K137").Value = 166767672
$ws.Range("L137").Value = 12445.3638
$ws.Range("M137").Value = -166762572
$ws.Range("N137").Value = -22645.3638

$ws.Range("H140").Value = 1825
$ws.Range("I140").Value = 1533.5714
$ws.Range("J140").Value = 2165
$ws.Range("K140").Value = 4600.7142
$ws.Range("L140").Value = 6495
$ws.Range("M140").Value = 579.2857999999997
$ws.Range("N140").Value = -16855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 800
$ws.Range("I13").Value = 800
$ws.Range("K13").Value = 800
$ws.Range("M13").Value = -661

$ws.Range("H126").Value = 2599.923
$ws.Range("I126").Value = 2367.2222
$ws.Range("J126").Value = 3123.5
$ws.Range("K126").Value = 7101.6666
$ws.Range("L126").Value = 9370.5
$ws.Range("M126").Value = -4631.6666
$ws.Range("N126").Value = -14310.5

$ws.Range("H132").Value = 2486.7273
$ws.Range("I132").Value = 2174.8235
$ws.Range("J132").Value = 3547.2
$ws.Range("K132").Value = 6524.470499999999
$ws.Range("L132").Value = 10641.6
$ws.Range("M132").Value = -3994.470499999999
$ws.Range("N132").Value = -15701.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3125
$ws.Range("I46").Value = 1666.6666
$ws.Range("J46").Value = 4583.3335
$ws.Range("K46").Value = 1666.6666
$ws.Range("L46").Value = 4583.3335
$ws.Range("M46").Value = -1478.6666
$ws.Range("N46").Value = -4959.3335

$ws.Range("H93").Value = 12170
$ws.Range("I93").Value = 15857.143
$ws.Range("K93").Value = 15857.143
$ws.Range("M93").Value = -14609.143

$ws.Range("H106").Value = 63333.332
$ws.Range("J106").Value = 63333.332
$ws.Range("L106").Value = 63333.332
$ws.Range("N106").Value = -65857.33199999999

$ws.Range("H132").Value = 2287.7827
$ws.Range("I132").Value = 2023.9333
$ws.Range("J132").Value = 2782.5
$ws.Range("K132").Value = 6071.7999
$ws.Range("L132").Value = 8347.5
$ws.Range("M132").Value = -3541.7999
$ws.Range("N132").Value = -13407.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3507
$ws.Range("I81").Value = 3485.2942
$ws.Range("J81").Value = 3599.25
$ws.Range("K81").Value = 6970.5884
$ws.Range("L81").Value = 7198.5
$ws.Range("M81").Value = -5909.5884
$ws.Range("N81").Value = -9320.5

$ws.Range("H84").Value = 3507
$ws.Range("I84").Value = 3485.2942
$ws.Range("J84").Value = 3599.25
$ws.Range("K84").Value = 34852.942
$ws.Range("L84").Value = 35992.5
$ws.Range("M84").Value = -29548.942
$ws.Range("N84").Value = -46600.5

$ws.Range("H117").Value = 24468
$ws.Range("J117").Value = 24468
$ws.Range("L117").Value = 24468
$ws.Range("N117").Value = -33646
